$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.526.95'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '1.811.20'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").Value = "'1.005"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").Value = "'304.80"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("D7").Value = "'0.4632"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.98%  '
$ws.Range("D8").Value = "'0.3561"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("D9").Value = "'46.16"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.55%  '
$ws.Range("D10").Value = "'0.07094"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("D11").Value = "'0.8985"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.19%  '
$ws.Range("D12").Value = "'0.07748"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").Value = "'19.31"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D14").Value = '1.822.03'
$ws.Range("E14").Value = '  +1.31%  '
$ws.Range("D15").Value = "'5.236"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = "'6.284"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.33%  '
$ws.Range("D17").Value = "'87.41"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.16%  '
$ws.Range("D18").Value = "'1.007"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").Value = "'0.000008521"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").Value = '26.576.88'
$ws.Range("E21").Value = '  +0.76%  '
$ws.Range("D22").Value = "'14.11"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = "'4.965"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").Value = "'10.52"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").Value = "'1.921"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.19%  '
$ws.Range("D26").Value = "'151.98"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").Value = "'17.81"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("D28").Value = "'2.005"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("D29").Value = "'112.57"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.39%  '
$ws.Range("D30").Value = "'4.790"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("D31").Value = "'0.08716"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.76%  '
$ws.Range("D32").Value = "'3.102"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.43%  '
$ws.Range("D33").Value = "'0.7296"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.22%  '
$ws.Range("D34").Value = "'4.416"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("D35").Value = "'2.704"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.04%  '
$ws.Range("D36").Value = "'1.116"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.62%  '
$ws.Range("D37").Value = "'1.070"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.85%  '
$ws.Range("D38").Value = "'0.01922"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.92%  '
$ws.Range("D39").Value = "'2.904"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.25%  '
$ws.Range("D40").Value = "'0.05061"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("D41").Value = "'0.5035"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.27%  '
$ws.Range("D42").Value = "'6.808"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.08%  '
$ws.Range("D43").Value = "'0.1492"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("D44").Value = "'7.956"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = "'0.4665"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = "'1.005"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("D47").Value = "'9.955"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.00%  '
$ws.Range("D48").Value = "'98.69"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("D49").Value = "'1.560"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("D50").Value = "'0.06011"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("D51").Value = "'63.42"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.02%  '
